# TC11_INS_CancerType-LiverCancer.xlsx edit
# Commit: "automation API poc changes & INS complete regression suite 23 scripts"
#
# The Programs-tab query in cell B2 is updated so the "Website" column is
# derived from prg.program_link / prg.program_acronym (via a CASE
# expression) instead of being a straight prg.website passthrough.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newQuery = @"
SELECT DISTINCT 
    prg.program_name AS "Program",
  CASE
    WHEN prg.program_link IS NOT NULL THEN prg.program_acronym
        ELSE prg.program_link
    END  AS "Website",
    prg.focus_area AS "Focus Area",
    prg.cancer_type AS "Cancer Type",
 CASE 
        WHEN prg.data_link IS NOT NULL THEN prg.website       
        ELSE prg.data_link
    END AS "Data Location Details"
FROM 
    df_program prg
WHERE 
     prg.cancer_type LIKE '%Liver Cancer%'
ORDER BY 
    lower(prg.program_name) ASC
LIMIT 100;
"@

$b2 = $ws.Range("B2")
$b2.Value = $newQuery
$b2.WrapText = $true
$b2.Font.Size = 12

# Re-select the cell the author ended up on after editing the query text.
$ws.Range("B8").Select()
